$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45171 -> 2023-09-02) for every
# data row (2 through 143). This update bumps that date by one day (45172 -> 2023-09-03)
# for all data rows, leaving the header (C1) untouched.
$ws.Range("C2:C143").Value = 45172
